$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix capitalization of existing headers: metadata4Ing -> metadata4ing
$ws.Range("D1").Value = "metadata4ing_IRI"
$ws.Range("E1").Value = "metadata4ing_DESC"

# Add new column F with the same header style as the other header cells
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "metadata4ing_DEF"

# Populate the new column F data rows
$ws.Range("F2").Value = "[locstr('A role is the function of an entity or agent with respect to an activity, in the context of a usage, generation, invalidation, association, start, and end.', 'en')]"
$ws.Range("F3").Value = "[]"
